{"js": "// Update the student name on the \"Name\" line.\nconst nameResults = context.document.body.search(\"B. Kanimozhi\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"P. Karthiga \", Word.InsertLocation.replace);\n}\n\n// Update the register number. In the source document this value is split\n// across two runs (\"Register No : 6114191040\" + \"29\"); searching/replacing\n// the full visible string collapses the result into a single run that keeps\n// the first run's formatting.\nconst regResults = context.document.body.search(\"Register No : 611419104029\", { matchCase: true });\nregResults.load(\"items\");\nawait context.sync();\n\nif (regResults.items.length > 0) {\n  regResults.items[0].insertText(\"Register No : 611419104031\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the student name on the \"Name\" line.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"B. Kanimozhi\"\n$find.Replacement.Text = \"P. Karthiga \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Update the register number. In the source document this value is split\n# across two runs (\"Register No : 6114191040\" + \"29\"); replacing the whole\n# visible string collapses it into a single run using the first run's\n# formatting.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Register No : 611419104029\"\n$find2.Replacement.Text = \"Register No : 611419104031\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
